$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 with new UUID/name/path
$ws.Range("A2").Value = "4d3c1167-0f83-461e-88a9-96df288c8e2b"
$ws.Range("B2").Value = "ueueuue ueueu ue"
$ws.Range("C2").Value = "src\renderer\storage\presentation-4d3c1167-0f83-461e-88a9-96df288c8e2b.xlsx"

# Add new row 3 with another record
$ws.Range("A3").Value = "106f270d-ab19-49df-9e5e-eb21a141fc58"
$ws.Range("B3").Value = "ii i ii isiad sii"
$ws.Range("C3").Value = "src\renderer\storage\presentation-106f270d-ab19-49df-9e5e-eb21a141fc58.xlsx"
